$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 49
$ws.Range("I2").Value = 156
$ws.Range("J2").Value = 624
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 193
$ws.Range("M2").Value = 12
$ws.Range("N2").Value = 116
$ws.Range("P2").Value = 1
$ws.Range("S2").Value = 70
$ws.Range("T2").Value = 101
$ws.Range("U2").Value = 11
$ws.Range("V2").Value = 937
$ws.Range("X2").Value = 974
$ws.Range("Y2").Value = 4
$ws.Range("Z2").Value = 13
$ws.Range("AA2").Value = 7
